$d = $word.ActiveDocument

function XmlEscape($s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

# Replaces the text of a single run (found via Find) with $newText, while leaving
# every other run/element in the paragraph (e.g. a preceding empty <w:r/>) untouched.
# $rPrXml is the literal <w:rPr>...</w:rPr> markup to keep on the run (or "" for none).
# $searchFrom is a Range to start searching from (defaults to the whole document);
# pass the previous match's range to continue past it for repeated text.
function Replace-RunText($oldText, $newText, $rPrXml, $searchFrom) {
    if ($searchFrom) {
        $rng = $d.Range($searchFrom.End, $d.Content.End)
    } else {
        $rng = $d.Content
    }
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $oldText)
        return $null
    }
    # Rebuild as a fresh Range; InsertXML on the Find-produced range object inserts
    # rather than replacing, but a Range freshly constructed over the same span replaces.
    $rng = $d.Range($rng.Start, $rng.End)

    $esc = XmlEscape $newText
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r>' + $rPrXml + '<w:t>' + $esc + '</w:t></w:r></w:p></w:body>' +
           '</w:document></pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xml)
    return $rng
}

# Title appears twice: Heading1 at the top (plain run), and a bold run near the end.
# Both map to the same new text.
$m1 = Replace-RunText "Play Bomber Squad Free: Unique 5-Reel Slot with Progressive Jackpot" "Play Bomber Squad | Free Slot Game" ""
Replace-RunText "Play Bomber Squad Free: Unique 5-Reel Slot with Progressive Jackpot" "Play Bomber Squad | Free Slot Game" "<w:rPr><w:b/></w:rPr>" $m1 | Out-Null

# "What we like" bullets
Replace-RunText "Thrilling storyline and well-designed symbols" "Cartoonish symbols drawn in detail" "" | Out-Null
Replace-RunText "Customizable gameplay experience" "Thrilling criminal storyline" "" | Out-Null
Replace-RunText "Five bonus options and progressive jackpot" "1024 ways to win" "" | Out-Null
Replace-RunText "1024 ways to win and potential for big payouts" "Customizable gameplay options" "" | Out-Null

# "What we don't like" bullets
Replace-RunText "RTP slightly lower than some other online slot games" "RTP value is within the average range" "" | Out-Null
Replace-RunText "Minimum bet may be too high for some players" "Limited bonus symbols" "" | Out-Null

# Meta description (italic paragraph)
Replace-RunText "Read a review of Bomber Squad, a five-reel online slot with 1024 ways to win, a thrilling storyline, and a unique bank robbery theme. Play for free now!" "Read our review of Bomber Squad and play this exciting slot game for free." "<w:rPr><w:i/></w:rPr>" | Out-Null
